$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '36.697.69'
$ws.Range("E2").Value = '  +4.06%  '
$ws.Range("D3").Value = '1.923.71'
$ws.Range("E3").Value = '  +2.47%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("E5").Value = '  +3.56%  '
$ws.Range("D6").Value = "'249.83"
$ws.Range("E6").Value = '  +1.33%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").Value = "'44.36"
$ws.Range("E8").Value = '  +1.63%  '
$ws.Range("D9").Value = "'59.05"
$ws.Range("E9").Value = '  +10.28%  '
$ws.Range("E10").Value = '  +3.92%  '
$ws.Range("E11").Value = '  +4.01%  '
$ws.Range("E12").Value = '  +2.52%  '
$ws.Range("D13").Value = "'14.60"
$ws.Range("E13").Value = '  +8.26%  '
$ws.Range("D14").Value = "'0.830"
$ws.Range("E14").Value = '  +8.51%  '
$ws.Range("D15").Value = '2.205.29'
$ws.Range("E15").Value = '  +2.55%  '
$ws.Range("E16").Value = '  +4.49%  '
$ws.Range("D17").Value = '1.922.42'
$ws.Range("E17").Value = '  +2.35%  '
$ws.Range("D18").Value = '36.675.56'
$ws.Range("E18").Value = '  +3.81%  '
$ws.Range("D19").Value = "'74.61"
$ws.Range("E19").Value = '  +2.82%  '
$ws.Range("D21").Value = "'251.63"
$ws.Range("E21").Value = '  +3.22%  '
$ws.Range("E22").Value = '  +4.66%  '
$ws.Range("D23").Value = "'5.27"
$ws.Range("E23").Value = '  +5.87%  '
$ws.Range("E24").Value = '  +2.06%  '
$ws.Range("E25").Value = '  -0.03%  '
$ws.Range("D26").Value = "'2.21"
$ws.Range("E26").Value = '  -1.39%  '
$ws.Range("D27").Value = "'168.17"
$ws.Range("E27").Value = '  +1.67%  '
$ws.Range("E28").Value = '  +3.73%  '
$ws.Range("E29").Value = '  +2.90%  '
$ws.Range("E30").Value = '  +2.17%  '
$ws.Range("D31").Value = "'4.60"
$ws.Range("E31").Value = '  +7.13%  '
$ws.Range("D32").Value = "'0.0621"
$ws.Range("E32").Value = '  +5.08%  '
$ws.Range("D33").Value = "'1.97"
$ws.Range("E33").Value = '  -3.70%  '
$ws.Range("E34").Value = '  +5.18%  '
$ws.Range("D35").Value = "'0.999"
$ws.Range("E35").Value = '  -0.14%  '
$ws.Range("D36").Value = "'1.54"
$ws.Range("E36").Value = '  -7.77%  '
$ws.Range("E37").Value = '  +17.79%  '
$ws.Range("D38").Value = "'0.911"
$ws.Range("E38").Value = '  +8.04%  '
$ws.Range("D39").Value = "'17.79"
$ws.Range("E39").Value = '  +49.89%  '
$ws.Range("D40").Value = "'2.07"
$ws.Range("E40").Value = '  +6.78%  '
$ws.Range("D41").Value = "'107.30"
$ws.Range("E41").Value = '  +11.45%  '
$ws.Range("D42").Value = "'0.0229"
$ws.Range("E42").Value = '  +5.01%  '
$ws.Range("D43").Value = "'17.30"
$ws.Range("E43").Value = '  -1.31%  '
$ws.Range("D44").Value = "'1.11"
$ws.Range("E44").Value = '  +3.78%  '
$ws.Range("D45").Value = "'2.74"
$ws.Range("E45").Value = '  +14.94%  '
$ws.Range("D46").Value = '1.344.13'
$ws.Range("E46").Value = '  +2.88%  '
$ws.Range("E47").Value = '  +1.26%  '
$ws.Range("D48").Value = "'0.0815"
$ws.Range("D49").Value = "'2.80"
$ws.Range("E49").Value = '  +3.20%  '
$ws.Range("D50").Value = "'6.46"
$ws.Range("E50").Value = '  +3.67%  '
$ws.Range("D51").Value = "'43.88"
$ws.Range("E51").Value = '  +4.28%  '
